# New weekly record is prepended at row 101 (the dataset is kept in
# reverse-chronological "most recent first" order within this block), so
# every existing record from row 101 down to row 247 shifts down by one
# row, and the sheet's dimension grows from A1:R247 to A1:R248.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101:247 down to 102:248, leaving a blank row 101 behind.
$ws.Rows.Item(101).Insert()

# The new row 101 shares every "template" column (Mercado ID, Mercado,
# Region, Codreg, Categoria ID, Categoria, Variedad, Calidad, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades, Clasificacion) with the rest of the series - only the
# date (D) and volume (J) are new for this entry. Grab that template from
# the row directly below (which now holds the data that used to live in
# row 101) and copy it up, then overwrite the two changed values.
$ws.Range("A102:R102").Copy()
$ws.Range("A101").PasteSpecial()

$ws.Cells.Item(101, 4).Value = 44579
$ws.Cells.Item(101, 10).Value = 160
